$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Add new values for the two new unit-label cells (D22, D24, D25).
#    Order matters for shared-string allocation: "m" must be registered
#    before "deg" so that the new shared strings come out as
#    index 48 = "m", index 49 = "deg" (matching the target workbook).
# ------------------------------------------------------------------
$ws.Range("D24").Value = "m"
$ws.Range("D25").Value = "m"
$ws.Range("D22").Value = "deg"

# ------------------------------------------------------------------
# Helper: apply a thin, automatic-colour box border around a range,
# reproducing the existing "bordered" cell style used throughout the
# sheet (e.g. the style already used by B2:F2, B3:H3, etc.).
# ------------------------------------------------------------------
function Add-ThinBorder($rng) {
    $borders = $rng.Borders
    $borders.LineStyle = 1      # xlContinuous
    $borders.Weight = 2         # xlThin
    $borders.ColorIndex = -4105 # xlAutomatic
}

# ------------------------------------------------------------------
# 2. Give the "Sideways on 30 degree slope" block its table-style
#    border (this is what turns style "3" cells into the bordered
#    style already used elsewhere in the sheet).
# ------------------------------------------------------------------
Add-ThinBorder $ws.Range("B20")
Add-ThinBorder $ws.Range("B21:C21")
Add-ThinBorder $ws.Range("B22:C22")
Add-ThinBorder $ws.Range("C23:D23")
Add-ThinBorder $ws.Range("C24")
Add-ThinBorder $ws.Range("C25")
Add-ThinBorder $ws.Range("B26:D26")
Add-ThinBorder $ws.Range("B27:D27")

# ------------------------------------------------------------------
# 3. Border the unit-label cells in column D (kg / deg / m / m).
# ------------------------------------------------------------------
Add-ThinBorder $ws.Range("D21")
Add-ThinBorder $ws.Range("D22")
Add-ThinBorder $ws.Range("D24")
Add-ThinBorder $ws.Range("D25")

# ------------------------------------------------------------------
# 4. Add borders around the bold row-labels (Mass/wheel base/weight
#    distance) that already use the bold label style.
# ------------------------------------------------------------------
Add-ThinBorder $ws.Range("B23:B25")

# ------------------------------------------------------------------
# 5. Update the view: scroll so row 6 is at the top and select
#    B20:D27 (activates B20 as the active cell, matching the target).
# ------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("B20:D27").Select()
